# decastro_r_assn04.xlsx -- "Add files via upload"
# Fill in BLAST hit rows 7-9 (previously blank placeholder rows), renumber the
# rank column, append a new row 12, widen column D, and update the selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Re-number the rank column A (was 2..8 across rows 5-11, becomes 1..7)
# ---------------------------------------------------------------------------
$ws.Range("A5").Value2  = 1
$ws.Range("A6").Value2  = 2
$ws.Range("A7").Value2  = 3
$ws.Range("A8").Value2  = 4
$ws.Range("A9").Value2  = 5
$ws.Range("A10").Value2 = 6
$ws.Range("A11").Value2 = 7

# ---------------------------------------------------------------------------
# 2. Row 6 gains a "Database description" value (column D)
# ---------------------------------------------------------------------------
$ws.Range("D6").Value = "watermelon_genes.fasta"

# ---------------------------------------------------------------------------
# 3. Rows 7-9 were empty placeholder rows; populate them with real BLAST hit
#    data. First copy the number/science formatting that already exists on
#    row 5 (G/H columns) and row 6 (E column) so no new styles are created,
#    then write the values/text.
# ---------------------------------------------------------------------------

# Donor single cells already carrying the exact formatting these cells need
# (single-cell sources tile cleanly across a multi-cell destination; a
# multi-cell source would instead repeat-tile past the destination bounds).
$ws.Range("G5").Copy()
$ws.Range("C7:C9").PasteSpecial(-4122)
$ws.Range("D7:D9").PasteSpecial(-4122)
$ws.Range("G7:G9").PasteSpecial(-4122)
$ws.Range("I9").PasteSpecial(-4122)

$ws.Range("H5").Copy()
$ws.Range("H7:H9").PasteSpecial(-4122)

$ws.Range("E6").Copy()
$ws.Range("E8:E9").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Row 7
$ws.Range("C7").Value2 = 303
$ws.Range("D7").Value = " watermelon_mt"
$ws.Range("E7").Value = "379,236 total"
$ws.Range("G7").Value2 = 560
$ws.Range("H7").Value2 = ("2E-161" -as [double])
$ws.Range("I7").Value = "plus/minus"

# Row 8
$ws.Range("C8").Value2 = 303
$ws.Range("D8").Value = "plant_mt"
$ws.Range("E8").Value2 = 11638289
$ws.Range("G8").Value2 = 560
$ws.Range("H8").Value2 = ("6E-160" -as [double])
$ws.Range("I8").Value = "plus/minus"

# Row 9
$ws.Range("C9").Value2 = 303
$ws.Range("D9").Value = "Nucleotide collection (nt)"
$ws.Range("E9").Value2 = 204700810597
$ws.Range("G9").Value2 = 560
$ws.Range("H9").Value2 = ("9E-156" -as [double])
$ws.Range("I9").Value = "lus/Minus"

# ---------------------------------------------------------------------------
# 4. New row 12, a copy of the same "Citrullus nad4L" placeholder row pattern
# ---------------------------------------------------------------------------
$ws.Range("A12").Value2 = 8
$ws.Range("B12").Value = "Citrullus nad4L"

# ---------------------------------------------------------------------------
# 5. Column D is a bit wider now that it holds longer descriptions, and no
#    longer needs to be a "best fit" column.
# ---------------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 26.6640625

# ---------------------------------------------------------------------------
# 6. Selection moves to C15
# ---------------------------------------------------------------------------
$ws.Range("C15").Select()
